$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validation")

$ws.Range("A3").Value = 1556
$ws.Range("B3").Value = 1856
$ws.Range("C3").Value = 2010
$ws.Range("D3").Value = 2515
$ws.Range("I3").Value = 93
$ws.Range("J3").Value = 101
$ws.Range("K3").Value = 108
$ws.Range("L3").Value = 122
